$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; G=77.63718533333333; H=232.911556; I=0.4513549673384918; J=0.4513549673384918; K=3; M=17.33599166666667; N=52.007975; O=0.4573561888773979; P=0.4573561888773979; Q=1345.917597962122; R=12113.2583816591; S=0.206429987692815; T=0.206429987692815 }
    3  = @{ E=3; G=77.63718533333333; H=232.911556; I=0.4513549673384918; J=0.4513549673384918; K=3; M=9.077707333333334; N=27.233122; O=0.2394870573052156; P=0.2394870573052156; Q=704.7676466397592; R=6342.908819757832; S=0.1080936729279871; T=0.1080936729279871 }
    4  = @{ E=3; G=77.63718533333333; H=232.911556; I=0.4513549673384918; J=0.4513549673384918; K=3; M=11.491094; N=34.473282; O=0.3031567538173866; P=0.3031567538173866; Q=892.1361945607546; R=8029.225751046791; S=0.1368313067176898; T=0.1368313067176897 }
    5  = @{ E=3; G=62.40815866666667; H=187.224476; I=0.3628188257432201; J=0.3628188257432201; K=3; M=17.33599166666667; N=52.007975; O=0.4573561888773979; P=0.4573561888773979; Q=1081.907318577345; R=9737.165867196101; S=0.1659374353948919; T=0.1659374353948919 }
    6  = @{ E=3; G=62.40815866666667; H=187.224476; I=0.3628188257432201; J=0.3628188257432201; K=3; M=9.077707333333334; N=27.233122; O=0.2394870573052156; P=0.2394870573052156; Q=566.5229995882304; R=5098.706996294072; S=0.08689041291217758; T=0.08689041291217758 }
    7  = @{ E=3; G=62.40815866666667; H=187.224476; I=0.3628188257432201; J=0.3628188257432201; K=3; M=11.491094; N=34.473282; O=0.3031567538173866; P=0.3031567538173866; Q=717.1380176055814; R=6454.242158450232; S=0.1099909774361507; T=0.1099909774361507 }
    8  = @{ E=3; G=31.96380833333333; H=95.891425; I=0.1858262069182881; J=0.1858262069182881; K=3; M=17.33599166666667; N=52.007975; O=0.4573561888773979; P=0.4573561888773979; Q=554.1243149015972; R=4987.118834114375; S=0.08498876578969099; T=0.08498876578969097 }
    9  = @{ E=3; G=31.96380833333333; H=95.891425; I=0.1858262069182881; J=0.1858262069182881; K=3; M=9.077707333333334; N=27.233122; O=0.2394870573052156; P=0.2394870573052156; Q=290.1580973087611; R=2611.42287577885; S=0.04450297146505091; T=0.0445029714650509 }
    10 = @{ E=3; G=31.96380833333333; H=95.891425; I=0.1858262069182881; J=0.1858262069182881; K=3; M=11.491094; N=34.473282; O=0.3031567538173866; P=0.3031567538173866; Q=367.2991261563166; R=3305.69213540685; S=0.0563344696635462; T=0.0563344696635462 }
}

foreach ($rowNum in $data.Keys) {
    $rowValues = $data[$rowNum]
    foreach ($col in $rowValues.Keys) {
        $ws.Range("$col$rowNum").Value = $rowValues[$col]
    }
}
